$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename move references from single-digit (z0bug.move_1 ...) to
# zero-padded two-digit (z0bug.move_01 ...) in columns A (id) and B
# (move_id), for rows 2 through 21 (moves 1-9; move_10 stays as-is).
for ($r = 2; $r -le 21; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)

    $valA = $cellA.Value2
    $valB = $cellB.Value2

    $matchA = $valA -match '^z0bug\.move_(\d)(_\d+)?$'
    if ($matchA) {
        $cellA.Value2 = "z0bug.move_0" + $matches[1] + $matches[2]
    }
    $matchB = $valB -match '^z0bug\.move_(\d)(_\d+)?$'
    if ($matchB) {
        $cellB.Value2 = "z0bug.move_0" + $matches[1] + $matches[2]
    }
}

# Selection moved to C5.
$null = $ws.Range("C5").Select()

# Column width tweaks (A and B slightly wider).
$ws.Columns.Item(1).ColumnWidth = 16.68
$ws.Columns.Item(2).ColumnWidth = 13.76
